$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the row containing "インテル株式会社" (row 59: #property=58)
$ws.Rows.Item(59).Delete()

# Renumber the "#property" serial column (A) for the rows that shifted up,
# so the sequence stays 1..157 with no gap.
for ($r = 59; $r -le 158; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
